$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Session Analysis Results")

$ws.Cells.Item(2, 7).Value = 'servinaz@med.asu.edu.eg, shaimaa.ahmed@med.asu.edu.eg, nourhan.mahmoud@med.asu.edu.eg, majorelle.magdy@med.asu.edu.eg, rana.abozaid@med.asu.edu.eg'
$ws.Cells.Item(10, 7).Value = 'Madeha.Saeed@med.asu.edu.eg, amany.raafat@med.asu.edu.eg, arwaelsayed03@med.asu.edu.eg, maryam.ahmed@med.asu.edu.eg, Marina.youhana@med.asu.edu.eg, esraa.mostafa@med.asu.edu.eg'
$ws.Cells.Item(18, 7).Value = 'shorokmohamed@med.asu.edu.eg, aya.hanafy@med.asu.edu.eg, Remon.Matta@med.asu.edu.eg, yasmin.m.senosy@med.asu.edu.eg'
$ws.Cells.Item(19, 7).Value = 'naema.gomaa@med.asu.edu.eg, Salma.hassan@med.asu.edu.eg, Monica.Eshak@med.asu.edu.eg, Remon.Matta@med.asu.edu.eg, yassmen.ahmed@med.asu.edu.eg, neveen.nashaat@med.asu.edu.eg, wafaa.ebida@med.asu.edu.eg, nardine.alfonse@med.asu.edu.eg'
$ws.Cells.Item(20, 7).Value = 'marinasorial@med.asu.edu.eg, Remon.Matta@med.asu.edu.eg, yassmen.ahmed@med.asu.edu.eg, youstina.magdy@med.asu.edu.eg, aya.hanafy@med.asu.edu.eg, nardine.alfonse@med.asu.edu.eg, Wafaa.ebida@med.asu.edu.eg'
$ws.Cells.Item(24, 7).Value = 'servinaz@med.asu.edu.eg, shaimaa.ahmed@med.asu.edu.eg, nourhan.mahmoud@med.asu.edu.eg, majorelle.magdy@med.asu.edu.eg, rana.abozaid@med.asu.edu.eg'
$ws.Cells.Item(25, 7).Value = 'manar.montaser@med.asu.edu.eg, backup@backdoor.com, alshimaa.atef@med.asu.edu.egm, gehanadel@med.asu.edu.eg'
$ws.Cells.Item(32, 7).Value = 'Madeha.Saeed@med.asu.edu.eg, amany.raafat@med.asu.edu.eg, arwaelsayed03@med.asu.edu.eg, maryam.ahmed@med.asu.edu.eg, Marina.youhana@med.asu.edu.eg, esraa.mostafa@med.asu.edu.eg'
$ws.Cells.Item(40, 7).Value = 'shorokmohamed@med.asu.edu.eg, aya.hanafy@med.asu.edu.eg, Remon.Matta@med.asu.edu.eg, yasmin.m.senosy@med.asu.edu.eg'
$ws.Cells.Item(41, 7).Value = 'naema.gomaa@med.asu.edu.eg, Salma.hassan@med.asu.edu.eg, Monica.Eshak@med.asu.edu.eg, Remon.Matta@med.asu.edu.eg, yassmen.ahmed@med.asu.edu.eg, neveen.nashaat@med.asu.edu.eg, wafaa.ebida@med.asu.edu.eg, nardine.alfonse@med.asu.edu.eg'
$ws.Cells.Item(42, 7).Value = 'marinasorial@med.asu.edu.eg, Remon.Matta@med.asu.edu.eg, yassmen.ahmed@med.asu.edu.eg, youstina.magdy@med.asu.edu.eg, aya.hanafy@med.asu.edu.eg, nardine.alfonse@med.asu.edu.eg, Wafaa.ebida@med.asu.edu.eg'
$ws.Cells.Item(46, 7).Value = 'shaimaa.ahmed@med.asu.edu.eg, nourhan.mahmoud@med.asu.edu.eg, hend_mahmoud@med.asu.edu.eg'
$ws.Cells.Item(52, 7).Value = 'mariam.noureldin@med.asu.edu.eg, Shimaa.ashraf@med.asu.edu.eg'
$ws.Cells.Item(54, 7).Value = 'basma.hamed@med.asu.edu.eg, Madeha.Saeed@med.asu.edu.eg, Eman.m.abosakaya@med.asu.edu.eg, marwa_mustafa@med.asu.edu.eg, merna.said@med.asu.edu.eg, yassmina.fattoh@med.asu.edu.eg, amany.raafat@med.asu.edu.eg, arwaelsayed03@med.asu.edu.eg, maryam.ahmed@med.asu.edu.eg, Marina.youhana@med.asu.edu.eg, maimustafa@med.asu.edu.eg'
$ws.Cells.Item(62, 7).Value = 'yassmen.ahmed@med.asu.edu.eg, shorokmohamed@med.asu.edu.eg, aya.hanafy@med.asu.edu.eg, wafaa.ebida@med.asu.edu.eg'
$ws.Cells.Item(64, 7).Value = 'wafaa.ebida@med.asu.edu.eg, youstina.magdy@med.asu.edu.eg'
$ws.Cells.Item(68, 7).Value = 'shaimaa.ahmed@med.asu.edu.eg, nourhan.mahmoud@med.asu.edu.eg, hend_mahmoud@med.asu.edu.eg'
$ws.Cells.Item(74, 7).Value = 'mariam.noureldin@med.asu.edu.eg, Shimaa.ashraf@med.asu.edu.eg'
$ws.Cells.Item(76, 7).Value = 'basma.hamed@med.asu.edu.eg, Madeha.Saeed@med.asu.edu.eg, Eman.m.abosakaya@med.asu.edu.eg, marwa_mustafa@med.asu.edu.eg, merna.said@med.asu.edu.eg, yassmina.fattoh@med.asu.edu.eg, amany.raafat@med.asu.edu.eg, arwaelsayed03@med.asu.edu.eg, maryam.ahmed@med.asu.edu.eg, Marina.youhana@med.asu.edu.eg, maimustafa@med.asu.edu.eg'
$ws.Cells.Item(84, 7).Value = 'yassmen.ahmed@med.asu.edu.eg, shorokmohamed@med.asu.edu.eg, aya.hanafy@med.asu.edu.eg, wafaa.ebida@med.asu.edu.eg'
$ws.Cells.Item(86, 7).Value = 'wafaa.ebida@med.asu.edu.eg, youstina.magdy@med.asu.edu.eg'
$ws.Cells.Item(90, 7).Value = 'Mohammedeltanany@med.asu.edu.eg, shaimaa.ahmed@med.asu.edu.eg, manar.montaser@med.asu.edu.eg'
$ws.Cells.Item(96, 7).Value = 'norhan.mohamed@med.asu.edu.eg, mariam.noureldin@med.asu.edu.eg, Sara_nabil@med.asu.edu.eg, aml.awwad@med.asu.edu.eg'
$ws.Cells.Item(98, 7).Value = 'basma.hamed@med.asu.edu.eg, Madeha.Saeed@med.asu.edu.eg, Eman.m.abosakaya@med.asu.edu.eg, marwa_mustafa@med.asu.edu.eg, merna.said@med.asu.edu.eg, yassmina.fattoh@med.asu.edu.eg, amany.raafat@med.asu.edu.eg, arwaelsayed03@med.asu.edu.eg, maryam.ahmed@med.asu.edu.eg, Marina.youhana@med.asu.edu.eg, maimustafa@med.asu.edu.eg'
$ws.Cells.Item(106, 7).Value = 'Monica.Eshak@med.asu.edu.eg, youstina.magdy@med.asu.edu.eg, neveen.nashaat@med.asu.edu.eg, wafaa.ebida@med.asu.edu.eg, nardine.alfonse@med.asu.edu.eg'
$ws.Cells.Item(107, 7).Value = 'Monica.Eshak@med.asu.edu.eg, yassmen.ahmed@med.asu.edu.eg, aya.hanafy@med.asu.edu.eg, neveen.nashaat@med.asu.edu.eg, wafaa.ebida@med.asu.edu.eg'
$ws.Cells.Item(108, 7).Value = 'marinasorial@med.asu.edu.eg, Remon.Matta@med.asu.edu.eg, yassmen.ahmed@med.asu.edu.eg, youstina.magdy@med.asu.edu.eg, aya.hanafy@med.asu.edu.eg, nardine.alfonse@med.asu.edu.eg, Wafaa.ebida@med.asu.edu.eg'
$ws.Cells.Item(112, 7).Value = 'Mohammedeltanany@med.asu.edu.eg, shaimaa.ahmed@med.asu.edu.eg, manar.montaser@med.asu.edu.eg'
$ws.Cells.Item(118, 7).Value = 'norhan.mohamed@med.asu.edu.eg, mariam.noureldin@med.asu.edu.eg, Sara_nabil@med.asu.edu.eg, aml.awwad@med.asu.edu.eg'
$ws.Cells.Item(120, 7).Value = 'basma.hamed@med.asu.edu.eg, Madeha.Saeed@med.asu.edu.eg, Eman.m.abosakaya@med.asu.edu.eg, marwa_mustafa@med.asu.edu.eg, merna.said@med.asu.edu.eg, yassmina.fattoh@med.asu.edu.eg, amany.raafat@med.asu.edu.eg, arwaelsayed03@med.asu.edu.eg, maryam.ahmed@med.asu.edu.eg, Marina.youhana@med.asu.edu.eg, maimustafa@med.asu.edu.eg'
$ws.Cells.Item(128, 7).Value = 'Monica.Eshak@med.asu.edu.eg, youstina.magdy@med.asu.edu.eg, neveen.nashaat@med.asu.edu.eg, wafaa.ebida@med.asu.edu.eg, nardine.alfonse@med.asu.edu.eg'
$ws.Cells.Item(129, 7).Value = 'Monica.Eshak@med.asu.edu.eg, yassmen.ahmed@med.asu.edu.eg, aya.hanafy@med.asu.edu.eg, neveen.nashaat@med.asu.edu.eg, wafaa.ebida@med.asu.edu.eg'
$ws.Cells.Item(130, 7).Value = 'marinasorial@med.asu.edu.eg, Remon.Matta@med.asu.edu.eg, yassmen.ahmed@med.asu.edu.eg, youstina.magdy@med.asu.edu.eg, aya.hanafy@med.asu.edu.eg, nardine.alfonse@med.asu.edu.eg, Wafaa.ebida@med.asu.edu.eg'
$ws.Cells.Item(134, 7).Value = 'Amira.Sobhy@med.asu.edu.eg, hend_mahmoud@med.asu.edu.eg, Veronia.rafat@med.asu.edu.eg, majorelle.magdy@med.asu.edu.eg, asmaa.reda@med.asu.edu.eg'
$ws.Cells.Item(142, 7).Value = 'basma.hamed@med.asu.edu.eg, marwa_mustafa@med.asu.edu.eg, merna.said@med.asu.edu.eg, yassmina.fattoh@med.asu.edu.eg, esraa.mostafa@med.asu.edu.eg'
$ws.Cells.Item(150, 7).Value = 'naema.gomaa@med.asu.edu.eg, Salma.hassan@med.asu.edu.eg, Monica.Eshak@med.asu.edu.eg, Remon.Matta@med.asu.edu.eg, yassmen.ahmed@med.asu.edu.eg, neveen.nashaat@med.asu.edu.eg, wafaa.ebida@med.asu.edu.eg, nardine.alfonse@med.asu.edu.eg'
$ws.Cells.Item(151, 7).Value = 'yassmen.ahmed@med.asu.edu.eg, wafaa.ebida@med.asu.edu.eg, marina_atef@med.asu.edu.eg, Monica.Eshak@med.asu.edu.eg'
$ws.Cells.Item(152, 7).Value = 'marina_atef@med.asu.edu.eg, Wafaa.ebida@med.asu.edu.eg'
$ws.Cells.Item(156, 7).Value = 'Mohammedeltanany@med.asu.edu.eg, mennatulla.medhat@med.asu.edu.eg, majorelle.magdy@med.asu.edu.eg, manar.montaser@med.asu.edu.eg, alshimaa.atef@med.asu.edu.egm'
$ws.Cells.Item(164, 7).Value = 'basma.hamed@med.asu.edu.eg, marwa_mustafa@med.asu.edu.eg, merna.said@med.asu.edu.eg, yassmina.fattoh@med.asu.edu.eg, esraa.mostafa@med.asu.edu.eg'
$ws.Cells.Item(165, 7).Value = 'Sarah.Mahdy@med.asu.edu.eg, Noran.Mahmoud@med.asu.edu.eg'
$ws.Cells.Item(172, 7).Value = 'naema.gomaa@med.asu.edu.eg, Salma.hassan@med.asu.edu.eg, Monica.Eshak@med.asu.edu.eg, Remon.Matta@med.asu.edu.eg, yassmen.ahmed@med.asu.edu.eg, neveen.nashaat@med.asu.edu.eg, wafaa.ebida@med.asu.edu.eg, nardine.alfonse@med.asu.edu.eg'
$ws.Cells.Item(173, 7).Value = 'yassmen.ahmed@med.asu.edu.eg, wafaa.ebida@med.asu.edu.eg, marina_atef@med.asu.edu.eg, Monica.Eshak@med.asu.edu.eg'
$ws.Cells.Item(174, 7).Value = 'marina_atef@med.asu.edu.eg, Wafaa.ebida@med.asu.edu.eg'
